# Update the 'F' column (想去人数 / interest counts) across all four
# worksheets to match the refreshed gh-pages data snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 1870
$ws.Range("F9").Value = 3645
$ws.Range("F10").Value = 1239
$ws.Range("F15").Value = 1520
$ws.Range("F17").Value = 1814
$ws.Range("F20").Value = 20
$ws.Range("F21").Value = 487
$ws.Range("F22").Value = 1571
$ws.Range("F26").Value = 1082
$ws.Range("F27").Value = 2341
$ws.Range("F28").Value = 396
$ws.Range("F29").Value = 7
$ws.Range("F30").Value = 4412
$ws.Range("F31").Value = 77
$ws.Range("F32").Value = 77
$ws.Range("F33").Value = 17

$ws = $wb.Worksheets.Item(2)
$ws.Range("F21").Value = 165
$ws.Range("F22").Value = 6
$ws.Range("F23").Value = 197
$ws.Range("F40").Value = 36
$ws.Range("F43").Value = 93

$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 2554
$ws.Range("F4").Value = 2563
$ws.Range("F6").Value = 167
$ws.Range("F10").Value = 3045
$ws.Range("F11").Value = 553
$ws.Range("F12").Value = 844
$ws.Range("F13").Value = 292
$ws.Range("F14").Value = 285

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 2554
$ws.Range("F3").Value = 167
$ws.Range("F8").Value = 553
$ws.Range("F9").Value = 844
$ws.Range("F10").Value = 1870
$ws.Range("F12").Value = 1239
$ws.Range("F14").Value = 1520
$ws.Range("F21").Value = 1814
$ws.Range("F23").Value = 20
$ws.Range("F24").Value = 487
$ws.Range("F26").Value = 1571
$ws.Range("F28").Value = 165
$ws.Range("F29").Value = 165
$ws.Range("F31").Value = 197
$ws.Range("F32").Value = 1082
$ws.Range("F34").Value = 2341
$ws.Range("F35").Value = 396
$ws.Range("F38").Value = 285
$ws.Range("F40").Value = 4412
$ws.Range("F41").Value = 77
$ws.Range("F43").Value = 17
$ws.Range("F47").Value = 93

